$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename headers ---------------------------------------------------
# "Offer Quantity *" -> "Quantity *"
$ws.Range("A1").Value = "Quantity *"
# "First Name *" -> "Full Name *" (will hold the merged first+last name)
$ws.Range("G1").Value = "Full Name *"

# --- Merge "First Name" (G) and "Last Name" (H) into a single column --
for ($r = 2; $r -le 7; $r++) {
    $first = $ws.Cells.Item($r, 7).Value()
    $last  = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 7).Value = "$first $last"
}

# --- Remove the now-redundant "Last Name *" column (H) -----------------
# Everything to the right (Address, Pan, Seller Signatory Emails, Bank
# Account, IFSC Code, Demat, City, Update Only, DP, Client) shifts one
# column to the left.
$ws.Columns("H:H").Delete()

# --- Repair hyperlinks that used to live in column K (now column J) ----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:emp3@myfirm.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:emp1@investor1.com")
$ws.Hyperlinks.Add($ws.Range("F7"), "mailto:emp1@investor2.com")
$ws.Hyperlinks.Add($ws.Range("J4"), "mailto:emp3@myfirm.com")
$ws.Hyperlinks.Add($ws.Range("J6"), "mailto:emp1@investor1.com")
$ws.Hyperlinks.Add($ws.Range("J7"), "mailto:emp1@investor2.com")

# --- Reflect the user's final selection on the merged column -----------
$ws.Range("G1:G7").Select()
